# Apply crypto price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.278.42"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "1.658.96"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.94%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.891.38"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.680.84"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.533"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.68%  "
$ws.Range("D17").Value = "27.270.38"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.98%  "
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").Value = "1.259.08"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  +4.27%  "
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.818"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").Value = "1.801.59"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0983"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("E51").Value = "  +0.39%  "
